$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(24800, 24800, 31, 11, 38, 6, 31, 11, 63, 6, $null, $null, $null, $null, $null, $null, $null, $null),
    @(23900, 23900, 31, 11, 63, 6, 37, 6, 63, 6, 37, 6, 38, 6, 37, 6, 39, 5),
    @(24200, 24200, 31, 11, 38, 6, 37, 6, 63, 6, $null, $null, $null, $null, $null, $null, $null, $null),
    @(24700, 24700, 31, 11, 63, 6, 37, 6, 63, 6, $null, $null, $null, $null, $null, $null, $null, $null),
    @(24600, 24600, 31, 11, 63, 6, 37, 6, 38, 6, $null, $null, $null, $null, $null, $null, $null, $null),
    @(24000, 24000, 31, 11, 38, 6, 37, 6, 63, 6, $null, $null, $null, $null, $null, $null, $null, $null),
    @(23800, 23800, 31, 11, 63, 6, 37, 6, 63, 6, $null, $null, $null, $null, $null, $null, $null, $null),
    @(24300, 24300, 31, 11, 63, 6, 37, 6, 63, 6, 31, 11, 38, 6, $null, $null, $null, $null)
)

# Columns: A B (skip C) D E F G (skip H) I J K L (skip M) N O P Q (skip R) S T U V
$cols = @(1, 2, 4, 5, 6, 7, 9, 10, 11, 12, 14, 15, 16, 17, 19, 20, 21, 22)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 3
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $val = $rowVals[$j]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $cols[$j]).Value = $val
        }
    }
}

$ws.Range("M15").Select()
